$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.495.13"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.854.85"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4748"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2760"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.90%  "
$ws.Range("D11").Value = "1.844.82"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.997"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6258"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").Value = "30.458.46"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "253.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +12.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9997"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007356"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.938"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.924"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "164.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.027"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.883"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1026"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.052"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.847"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04848"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.134"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7013"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.691"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01898"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8766"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.993"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9995"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4080"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.531"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.197"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1204"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.543"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05497"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.356"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3705"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.33%  "
